# Generate Report for Handback
# Updates the localization-status workbook to reflect that the de-de
# handback has completed (zh-cn had already completed earlier in this run).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusHandedBack = "Handed back: in sync with en-US"
$targetMdUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6146d72ab9485d28aa02533a25a550baaf089de8/e2e/2a6ce15f-50ee-4d58-ab78-d69324521b6e.md"
$targetMdDisplay  = "2a6ce15f-50ee-4d58-ab78-d69324521b6e.md"
$zhCnXlfName      = "2a6ce15f-50ee-4d58-ab78-d69324521b6e.da69cf47134c556097c2f17e99a70b816d809605.zh-cn.xlf"
$deDeXlfName      = "2a6ce15f-50ee-4d58-ab78-d69324521b6e.da69cf47134c556097c2f17e99a70b816d809605.de-de.xlf"
$zhCnHandbackTime = "2016-08-29 17:10:53"
$deDeHandbackTime = "2016-08-29 17:11:01"

# ---------------------------------------------------------------------------
# 1. Status column now reads "Handed back: in sync with en-US" everywhere it
#    used to read "Ready for handoff" (Overview + both language sheets).
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack

$wsZhCn.Range("C2").Value = $statusHandedBack
$wsZhCn.Range("C3").Value = $statusHandedBack

$wsDeDe.Range("C2").Value = $statusHandedBack
$wsDeDe.Range("C3").Value = $statusHandedBack

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: the handback run for this language finished earlier, so the
#    "Latest Target File" / "Latest Handback File" / "Latest Handback
#    DateTime" columns are now populated for both rows.
# ---------------------------------------------------------------------------
$wsZhCn.Range("I2").Value = $targetMdDisplay
$wsZhCn.Range("J2").Value = $zhCnXlfName
$wsZhCn.Range("K2").Value = $zhCnHandbackTime

$wsZhCn.Range("I3").Value = $targetMdDisplay
$wsZhCn.Range("J3").Value = $zhCnXlfName
$wsZhCn.Range("K3").Value = $zhCnHandbackTime

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $targetMdUrl, "", "", $targetMdDisplay) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $targetMdUrl, "", "", $targetMdDisplay) | Out-Null

$wsZhCn.Range("I2").Font.Color = 15570276
$wsZhCn.Range("I2").Font.Underline = $true
$wsZhCn.Range("I2").Font.Name = "Calibri"
$wsZhCn.Range("I3").Font.Color = 15570276
$wsZhCn.Range("I3").Font.Underline = $true
$wsZhCn.Range("I3").Font.Name = "Calibri"

# ---------------------------------------------------------------------------
# 3. de-de sheet: the handback run for this language just completed, so the
#    same three columns are populated with the de-de artifacts / timestamp.
# ---------------------------------------------------------------------------
$wsDeDe.Range("I2").Value = $targetMdDisplay
$wsDeDe.Range("J2").Value = $deDeXlfName
$wsDeDe.Range("K2").Value = $deDeHandbackTime

$wsDeDe.Range("I3").Value = $targetMdDisplay
$wsDeDe.Range("J3").Value = $deDeXlfName
$wsDeDe.Range("K3").Value = $deDeHandbackTime

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $targetMdUrl, "", "", $targetMdDisplay) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $targetMdUrl, "", "", $targetMdDisplay) | Out-Null

$wsDeDe.Range("I2").Font.Color = 15570276
$wsDeDe.Range("I2").Font.Underline = $true
$wsDeDe.Range("I2").Font.Name = "Calibri"
$wsDeDe.Range("I3").Font.Color = 15570276
$wsDeDe.Range("I3").Font.Underline = $true
$wsDeDe.Range("I3").Font.Name = "Calibri"

# ---------------------------------------------------------------------------
# 4. Widen the columns that now hold longer content: the Status column (now
#    the long "Handed back: in sync with en-US" string) and the two new
#    filename columns (Latest Target File / Latest Handback File).
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.17   # E: zh-cn status
$wsOverview.Columns.Item(6).ColumnWidth = 29.17   # F: de-de status

$wsZhCn.Columns.Item(3).ColumnWidth = 29.17        # C: Status
$wsZhCn.Columns.Item(9).ColumnWidth = 39.2         # I: Latest Target File
$wsZhCn.Columns.Item(10).ColumnWidth = 39.2        # J: Latest Handback File

$wsDeDe.Columns.Item(3).ColumnWidth = 29.17        # C: Status
$wsDeDe.Columns.Item(9).ColumnWidth = 39.2         # I: Latest Target File
$wsDeDe.Columns.Item(10).ColumnWidth = 39.2        # J: Latest Handback File

Write-Output "done"
